$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" -------------------
# Overview sheet: per-language status cells (zh-cn column E, de-de column F)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn detail sheet: Status column (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de detail sheet: Status column (C)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width shrink for the Status columns now that the text is shorter
# Original autofit width was ~17.22 chars ("Ready for handoff"); the new
# text ("In Translation") is narrower, so the recalculated best-fit width for
# every "Status" column (Overview!E:F, zh-cn!C, de-de!C) is ~13.41 chars.
# (ColumnWidth is requested in character units; 12.5 is the input that lands
# the host's internal width snapping nearest to that ~13.41 target.)
$newStatusWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
